$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-20 Sunday" "2024-10-21 Monday"

Replace-Text "940×6=" "780×9="
Replace-Text "915×7=" "725×8="
Replace-Text "979×9=" "246×6="
Replace-Text "779×5=" "631×7="
Replace-Text "904×8=" "554×9="

Replace-Text "705×6=" "860×2="
Replace-Text "346×7=" "136×6="
Replace-Text "697×7=" "917×9="
Replace-Text "125×2=" "404×7="
Replace-Text "268×4=" "942×8="

Replace-Text "728×9=" "886×8="
Replace-Text "486×2=" "243×2="
Replace-Text "977×6=" "313×7="
Replace-Text "913×4=" "167×8="
Replace-Text "147×2=" "797×3="

Replace-Text "181×7=" "362×9="
Replace-Text "151×2=" "431×3="
Replace-Text "260×4=" "447×5="
Replace-Text "990×4=" "499×6="
Replace-Text "893×7=" "119×3="

Replace-Text "410×2=" "970×9="
Replace-Text "317×6=" "478×4="
Replace-Text "257×8=" "679×8="
Replace-Text "521×7=" "804×2="
Replace-Text "220×8=" "463×3="
